$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- H3: tag the "Code start" milestone column with a new note ------------
$ws.Range("H3").Value = "Basic functionality"

# --- New "X" marks (fill style) on I8 and G16, copied from an existing ----
# --- marked cell so the style index (s="2") is reused rather than a new ---
# --- style being created. -------------------------------------------------
$ws.Range("F8").Copy()
$ws.Range("I8").PasteSpecial(-4122)
$ws.Range("H16").Copy()
$ws.Range("G16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Sprint effort totals (row 26) were retyped during the meeting --------
$ws.Range("F26").Value = 12
$ws.Range("G26").Value = 16
$ws.Range("H26").Value = 16
$ws.Range("I26").Value = 17

# --- Drop the old "average per sprint" helper row (B26/4) -----------------
$ws.Range("B27").ClearContents()

# --- Widen the Sprint0/1/2 columns slightly (bestFit columns C:E) ---------
# (target stored width is 10.6328125 "characters"; the host's ColumnWidth
# setter quantises to whole pixels, so 9.8 is the input that lands closest
# to that stored value after the pixel round-trip.)
$ws.Range("C1:E1").EntireColumn.ColumnWidth = 9.8

# --- View state: zoomed in to 160%, cursor moved back up to F9 ------------
$excel.ActiveWindow.Zoom = 160
$ws.Range("F9").Select() | Out-Null
